# Add a new service event row (row 19) to the Card20 sheet and fill the
# previously-empty "nan" placeholder cells in row 18.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card20")

# --- Row 18: the blank inline-string cells get the literal text "nan" ---
$ws.Range("B18").Value = "nan"
$ws.Range("C18").Value = "nan"
$ws.Range("D18").Value = "nan"
$ws.Range("E18").Value = "nan"
$ws.Range("F18").Value = "nan"
$ws.Range("G18").Value = "nan"
$ws.Range("H18").Value = "nan"
$ws.Range("I18").Value = "nan"
$ws.Range("J18").Value = "nan"
$ws.Range("K18").Value = "nan"
$ws.Range("M18").Value = "nan"

# --- Row 19: brand-new service event ---
# "20" must stay TEXT (matches the rest of column A) rather than being
# auto-coerced to a number; force text via NumberFormat, then restore the
# default "Normal" style so no stray number-format sticks to the cell.
$ws.Range("A19").NumberFormat = "@"
$ws.Range("A19").Value = "20"
$ws.Range("A19").Style = "Normal"

$ws.Range("L19").Value = "16\12\2024"
$ws.Range("N19").Value = "تم عمل صيانه نصف سنويه"
$ws.Range("O19").Value = "تيم العمل"
